# Regenerate save_data to use K instead of Strike#, writing updated s_vals (K column, col G)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K (col G) values computed from recalculated std/mean, per row
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 1
